{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Summary of the edit (per the canonical OOXML diff):\n//  - Paragraph 1: adjacent runs that only differed because of\n//    <w:proofErr> spell-check markers are merged back into single runs\n//    (\"Hello \" + \"World\" -> \"Hello World\", \"Foo\" + \" Bar\" -> \"Foo Bar\",\n//    the run holding <w:tab/> absorbs \"Github is the best\", and\n//    \" \" + \"asdasd\" -> \" asdasd\"); all <w:proofErr> elements are gone.\n//  - Paragraph 2: same kind of run-merge (\"A\" + \"sd\" stay separate runs,\n//    but \" \" + \"fsdfsadfa\" -> \" fsdfsadfa\"); <w:proofErr> elements removed.\n//  - Paragraph 3: unchanged empty paragraph.\n//  - Paragraph 4: the other empty paragraph is removed (merged away).\n//  - Paragraph 5 (becomes paragraph 4): all its runs collapse into a\n//    single run with fresh text \"aaaaaaaaaa\", and it loses its empty\n//    <w:pPr> (no explicit paragraph mark run properties any more).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Helper: wrap a <w:body> inner fragment into a full OOXML package so it\n// can be fed to Range.insertOoxml() for exact, deterministic markup.\nfunction wrapBodyOoxml(innerXml) {\n  return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n    '<pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + innerXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n}\n\nconst boldRpr = '<w:rPr><w:b/><w:bCs/></w:rPr>';\n\n// --- Paragraph 1: \"Hello World\" / \"Foo Bar\" / tab+\"Github is the best\" /\n//     \" asdasd\" / \"asd\" ---------------------------------------------------\nconst para1Xml =\n  '<w:p><w:pPr>' + boldRpr + '</w:pPr>' +\n  '<w:r><w:t>Hello World</w:t></w:r>' +\n  '<w:r>' + boldRpr + '<w:t>Foo Bar</w:t></w:r>' +\n  '<w:r>' + boldRpr + '<w:tab/><w:t>Github is the best</w:t></w:r>' +\n  '<w:r>' + boldRpr + '<w:t xml:space=\"preserve\"> asdasd</w:t></w:r>' +\n  '<w:r>' + boldRpr + '<w:t>asd</w:t></w:r>' +\n  '</w:p>';\nparagraphs.items[0].getRange().insertOoxml(wrapBodyOoxml(para1Xml), \"Replace\");\nawait context.sync();\n\n// --- Paragraph 2: \"A\" / \"sd\" / \" fsdfsadfa\" -----------------------------\nconst para2Xml =\n  '<w:p><w:pPr>' + boldRpr + '</w:pPr>' +\n  '<w:r>' + boldRpr + '<w:t>A</w:t></w:r>' +\n  '<w:r>' + boldRpr + '<w:t>sd</w:t></w:r>' +\n  '<w:r>' + boldRpr + '<w:t xml:space=\"preserve\"> fsdfsadfa</w:t></w:r>' +\n  '</w:p>';\nparagraphs.items[1].getRange().insertOoxml(wrapBodyOoxml(para2Xml), \"Replace\");\nawait context.sync();\n\n// --- Paragraph 3: leave untouched (still an empty bold paragraph) ------\n\n// --- Paragraph 4: delete outright; its content merges away, leaving the\n//     previous paragraph 5's content as the new last paragraph ---------\nparagraphs.items[3].delete();\nawait context.sync();\n\n// Re-fetch paragraphs after the delete shifted indices.\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Former paragraph 5 (now index 3): single run \"aaaaaaaaaa\", and the\n//     paragraph mark no longer carries explicit bold run properties ----\nconst para5Xml = '<w:p><w:r>' + boldRpr + '<w:t>aaaaaaaaaa</w:t></w:r></w:p>';\nparagraphs.items[3].getRange().insertOoxml(wrapBodyOoxml(para5Xml), \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document (commonly aliased $d below).\n#\n# Summary of the edit (per the canonical OOXML diff):\n#  - Paragraph 1: adjacent runs that only differed because of\n#    <w:proofErr> spell-check markers are merged back into single runs\n#    (\"Hello \" + \"World\" -> \"Hello World\", \"Foo\" + \" Bar\" -> \"Foo Bar\",\n#    the run holding <w:tab/> absorbs \"Github is the best\", and\n#    \" \" + \"asdasd\" -> \" asdasd\"); all <w:proofErr> elements are gone.\n#  - Paragraph 2: same kind of run-merge (\"A\" + \"sd\" stay separate runs,\n#    but \" \" + \"fsdfsadfa\" -> \" fsdfsadfa\"); <w:proofErr> elements removed.\n#  - Paragraph 3: unchanged empty paragraph.\n#  - Paragraph 4: the other empty paragraph is removed (merged away).\n#  - Paragraph 5 (becomes paragraph 4): all its runs collapse into a\n#    single run with fresh text \"aaaaaaaaaa\", and it loses its empty\n#    <w:pPr> (no explicit paragraph mark run properties any more).\n\n$d = $word.ActiveDocument\n\n$boldRpr = \"<w:rPr><w:b/><w:bCs/></w:rPr>\"\n\nfunction Wrap-BodyOoxml($innerXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n        '<pkg:xmlData>' +\n        '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n        '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n        '</Relationships>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $innerXml + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\n# --- Paragraph 1: \"Hello World\" / \"Foo Bar\" / tab+\"Github is the best\" /\n#     \" asdasd\" / \"asd\" ---------------------------------------------------\n$para1Xml = '<w:p><w:pPr>' + $boldRpr + '</w:pPr>' +\n    '<w:r><w:t>Hello World</w:t></w:r>' +\n    '<w:r>' + $boldRpr + '<w:t>Foo Bar</w:t></w:r>' +\n    '<w:r>' + $boldRpr + '<w:tab/><w:t>Github is the best</w:t></w:r>' +\n    '<w:r>' + $boldRpr + '<w:t xml:space=\"preserve\"> asdasd</w:t></w:r>' +\n    '<w:r>' + $boldRpr + '<w:t>asd</w:t></w:r>' +\n    '</w:p>'\n$d.Paragraphs.Item(1).Range.InsertXML((Wrap-BodyOoxml $para1Xml))\n\n# --- Paragraph 2: \"A\" / \"sd\" / \" fsdfsadfa\" -----------------------------\n$para2Xml = '<w:p><w:pPr>' + $boldRpr + '</w:pPr>' +\n    '<w:r>' + $boldRpr + '<w:t>A</w:t></w:r>' +\n    '<w:r>' + $boldRpr + '<w:t>sd</w:t></w:r>' +\n    '<w:r>' + $boldRpr + '<w:t xml:space=\"preserve\"> fsdfsadfa</w:t></w:r>' +\n    '</w:p>'\n$d.Paragraphs.Item(2).Range.InsertXML((Wrap-BodyOoxml $para2Xml))\n\n# --- Paragraph 3: leave untouched (still an empty bold paragraph) ------\n\n# --- Paragraph 4: delete outright; its content merges away, leaving the\n#     previous paragraph 5's content as the new last paragraph ---------\n$d.Paragraphs.Item(4).Range.Delete()\n\n# --- Former paragraph 5 (now paragraph 4): single run \"aaaaaaaaaa\", and\n#     the paragraph mark no longer carries explicit bold run properties -\n$para5Xml = '<w:p><w:r>' + $boldRpr + '<w:t>aaaaaaaaaa</w:t></w:r></w:p>'\n$d.Paragraphs.Item(4).Range.InsertXML((Wrap-BodyOoxml $para5Xml))\n"}
